$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.768.18'
$ws.Range("E2").Value = '  +0.90%  '
$ws.Range("D3").Value = '1.649.17'
$ws.Range("E3").Value = '  +1.30%  '
$ws.Range("E4").Value = '  +0.55%  '
$ws.Range("D5").Value = '''216.61'
$ws.Range("E5").Value = '  +1.53%  '
$ws.Range("E6").Value = '  +0.37%  '
$ws.Range("E7").Value = '  +0.45%  '
$ws.Range("E8").Value = '  +0.90%  '
$ws.Range("E9").Value = '  +0.66%  '
$ws.Range("D10").Value = '''19.24'
$ws.Range("E10").Value = '  +2.42%  '
$ws.Range("E11").Value = '  -0.01%  '
$ws.Range("D12").Value = '1.877.83'
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.643.36'
$ws.Range("E13").Value = '  +2.52%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '''4.20'
$ws.Range("E14").Value = '  +1.36%  '
$ws.Range("D15").Value = '''0.533'
$ws.Range("E15").Value = '  +1.76%  '
$ws.Range("D16").Value = '''65.38'
$ws.Range("E16").Value = '  +0.46%  '
$ws.Range("D17").Value = '26.779.22'
$ws.Range("E17").Value = '  +0.87%  '
$ws.Range("D18").Value = '0.0₃0744'
$ws.Range("E18").Value = '  +0.44%  '
$ws.Range("D19").Value = '''217.70'
$ws.Range("E19").Value = '  +1.37%  '
$ws.Range("E20").Value = '  +0.45%  '
$ws.Range("E21").Value = '  +1.95%  '
$ws.Range("D22").Value = '''2.47'
$ws.Range("E22").Value = '  +15.38%  '
$ws.Range("D23").Value = '''6.27'
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("D24").Value = '''9.49'
$ws.Range("E24").Value = '  +1.78%  '
$ws.Range("D25").Value = '''147.20'
$ws.Range("E25").Value = '  -0.13%  '
$ws.Range("E26").Value = '  +0.50%  '
$ws.Range("E27").Value = '  -0.23%  '
$ws.Range("D28").Value = '''7.18'
$ws.Range("E28").Value = '  +4.04%  '
$ws.Range("D29").Value = '''15.78'
$ws.Range("E29").Value = '  +1.39%  '
$ws.Range("E30").Value = '  +1.28%  '
$ws.Range("E31").Value = '  +1.51%  '
$ws.Range("D32").Value = '''3.36'
$ws.Range("E33").Value = '  +1.56%  '
$ws.Range("D34").Value = '1.281.48'
$ws.Range("E34").Value = '  +3.13%  '
$ws.Range("D35").Value = '''1.55'
$ws.Range("E35").Value = '  +2.99%  '
$ws.Range("E36").Value = '  +2.88%  '
$ws.Range("E37").Value = '  +2.10%  '
$ws.Range("E38").Value = '  +5.78%  '
$ws.Range("E39").Value = '  +4.41%  '
$ws.Range("E40").Value = '  +0.51%  '
$ws.Range("E41").Value = '  +2.22%  '
$ws.Range("E42").Value = '  -0.82%  '
$ws.Range("E43").Value = '  +2.10%  '
$ws.Range("D44").Value = '1.788.92'
$ws.Range("E44").Value = '  +1.39%  '
$ws.Range("D45").Value = '''92.02'
$ws.Range("E45").Value = '  -1.33%  '
$ws.Range("D46").Value = '''59.80'
$ws.Range("E46").Value = '  +9.02%  '
$ws.Range("D47").Value = '''1.61'
$ws.Range("E47").Value = '  +1.37%  '
$ws.Range("E48").Value = '  -0.34%  '
$ws.Range("D49").Value = '''0.0515'
$ws.Range("E49").Value = '  +1.10%  '
$ws.Range("D50").Value = '''7.76'
$ws.Range("E50").Value = '  +3.45%  '
$ws.Range("D51").Value = '''0.0975'
$ws.Range("E51").Value = '  +1.75%  '
